$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "41.240.08"
$c.ClearFormats()
$ws.Cells.Item(2, 5).Value = "  -1.85%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "2.178.77"
$c.ClearFormats()
$ws.Cells.Item(3, 5).Value = "  -1.99%  "
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "237.29"
$c.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -2.38%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.614"
$c.ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -0.98%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "70.43"
$c.ClearFormats()
$ws.Cells.Item(7, 5).Value = "  -5.34%  "
$ws.Cells.Item(8, 5).Value = "  +0.05%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.580"
$c.ClearFormats()
$ws.Cells.Item(9, 5).Value = "  -6.26%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "40.23"
$c.ClearFormats()
$ws.Cells.Item(10, 5).Value = "  -9.24%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.0930"
$c.ClearFormats()
$ws.Cells.Item(11, 5).Value = "  -3.58%  "
$ws.Cells.Item(12, 5).Value = "  -2.40%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "6.77"
$c.ClearFormats()
$ws.Cells.Item(13, 5).Value = "  -5.81%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "2.501.97"
$c.ClearFormats()
$ws.Cells.Item(14, 5).Value = "  -2.05%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "13.94"
$c.ClearFormats()
$ws.Cells.Item(15, 5).Value = "  -2.64%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.811"
$c.ClearFormats()
$ws.Cells.Item(16, 5).Value = "  -4.33%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "2.183.03"
$c.ClearFormats()
$ws.Cells.Item(17, 5).Value = "  -1.20%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "41.039.75"
$c.ClearFormats()
$ws.Cells.Item(18, 5).Value = "  -2.18%  "
$ws.Cells.Item(19, 5).Value = "  -7.88%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "70.51"
$c.ClearFormats()
$ws.Cells.Item(20, 5).Value = "  -2.83%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "5.96"
$c.ClearFormats()
$ws.Cells.Item(21, 5).Value = "  -4.11%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "10.10"
$c.ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -10.60%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "225.91"
$c.ClearFormats()
$ws.Cells.Item(23, 5).Value = "  -1.91%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "1.95"
$c.ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -6.95%  "
$ws.Cells.Item(25, 5).Value = "  +0.18%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "10.88"
$c.ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -6.40%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "3.55"
$c.ClearFormats()
$ws.Cells.Item(27, 5).Value = "  -1.25%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "2.22"
$c.ClearFormats()
$ws.Cells.Item(28, 5).Value = "  -2.93%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "167.10"
$c.ClearFormats()
$ws.Cells.Item(30, 5).Value = "  +0.18%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "20.01"
$c.ClearFormats()
$ws.Cells.Item(31, 5).Value = "  -3.09%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "31.21"
$c.ClearFormats()
$ws.Cells.Item(32, 5).Value = "  +5.69%  "
$ws.Cells.Item(33, 5).Value = "  -4.36%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "5.18"
$c.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -9.96%  "
$ws.Cells.Item(35, 5).Value = "  -3.13%  "
$ws.Cells.Item(36, 5).Value = "  -9.44%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "4.13"
$c.ClearFormats()
$ws.Cells.Item(37, 5).Value = "  -4.01%  "
$ws.Cells.Item(38, 5).Value = "  -5.67%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "12.24"
$c.ClearFormats()
$ws.Cells.Item(39, 5).Value = "  -6.04%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "2.08"
$c.ClearFormats()
$ws.Cells.Item(40, 5).Value = "  -3.02%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "5.43"
$c.ClearFormats()
$ws.Cells.Item(41, 5).Value = "  -4.23%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "60.25"
$c.ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -7.51%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.191"
$c.ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -4.63%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "8.34"
$c.ClearFormats()
$ws.Cells.Item(44, 5).Value = "  -4.93%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.0973"
$c.ClearFormats()
$ws.Cells.Item(45, 5).Value = "  -3.78%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "98.33"
$c.ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -5.99%  "
$ws.Cells.Item(47, 5).Value = "  -2.85%  "
$ws.Cells.Item(48, 5).Value = "  -2.91%  "
$ws.Cells.Item(49, 5).Value = "  -8.63%  "
$ws.Cells.Item(50, 5).Value = "  -2.85%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "2.379.00"
$c.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -2.04%  "

Write-Host "Updated crypto prices and volume percentages."
